$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

# --- Row 2 ---
$ws.Range("A2").Value = "281474991395097-1749578321343"
$ws.Range("B2").Value = "Harsh Brake"
$ws.Range("C2").Value = "2025-06-10T11:58:41.343"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "281474991395097"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "125"
$ws.Range("E2").ClearFormats()

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "51834055"
$ws.Range("F2").ClearFormats()

$ws.Range("G2").Value = "DAVID SERRANO"
$ws.Range("H2").Value = 20.68131202
$ws.Range("I2").Value = -103.38573738
$ws.Range("J2").Value = 0.7485191226005554
$ws.Range("K2").Value = "No video URL"
$ws.Range("L2").Value = "No video URL"

# --- Row 3 ---
$ws.Range("A3").Value = "281474991205821-1749578080783"
$ws.Range("B3").Value = "No Seat Belt"
$ws.Range("C3").Value = "2025-06-10T11:54:40.783"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "281474991205821"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "148"
$ws.Range("E3").ClearFormats()

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "51834015"
$ws.Range("F3").ClearFormats()

$ws.Range("G3").Value = "LUIS IBARRA"
$ws.Range("H3").Value = 20.56918449
$ws.Range("I3").Value = -103.45711918
$ws.Range("K3").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1749578078283/gFqVKuRiGj-camera-video-segment-driver-1749578080783.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOSQKCUFZ%2F20250611%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250611T180153Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEAIaCXVzLXdlc3QtMiJIMEYCIQDhk5Onjj%2Fjx%2BlzJ8o205%2Bv1tsdlBBdthpD%2F4QpzF6bYwIhALNEtuuSmacwYjQvcc0fMlCNBKXYo4J2gNKCXX%2B%2BFrywKuYDCNv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgzUz0HnBnoKLwdoez8qugO3MhgIwL3LclI71xVisID0mBDNs8067PMeYBZe%2Fko42MpEUpV0qkdrO1tB%2B44kCAwuoywLIdW7BVZn1dRlGZoEle6YSkylm5ods18ww5tyy1Opr0i9RNZo%2BWs433CKDcDpNBVuGQ%2FGLm9iJJSabX8z6Ncao2cxw5RY7TBfDD5rs5JyPkyZLqDrqyL2wESG2b%2BMSQDTjjuv5%2B9dZKJlmlS7IbzOKAV7i3L4Nq3oRxfEhqClmurbIC0HYYb%2FIpqOVoeR7snITMw3Fs2FKMJuZogmbEwxShPXmf%2FOT7pD62oVaAXztypnUSALnv9ijYdXa9khUkanDDiiHoqZh7EdHrh9XwYQwi6wHj03C%2FRsQCWs2e3IyKkc0KwPpBnW6qRa3%2FSZ5Wr5Vj2WUKkSD9apEfa%2FxAwj5aTpn3h%2FdkLWTdI5s1ajKq3FXIcogbcrtGHxSnPXo9JKDh53iscfm7ltH%2Bf5E24SkIcBuxL0oDtgUmZPzU8%2FqpGPF3TOi0Vh05rzhZ6hyzR62J7%2Bhnd4PzqLIMjNlJqhjlD6lGlgW4SRix1egzOpUWWPZ6g05OY7J8TcYqnvh%2BzwGPNzIA1BMLeAp8IGOqQBbYFp2vYp3UIZhng%2B7s9p23ghBnFOaP8tnV%2FZYfvpObScfaqtZixTQfb4CohqRC5nEZNktkMpWTNjEiBcvAyB6bYR8b4fiLQGFpaY0xMTYZZhp3%2FrIEyLSyD8Tl0XWURGYalzMic%2Bf4IYUvGCh5lTYLJ3%2FX8ltWZzshJtOT0LY80f0Q6mYJoFgv%2F7677VGhbhlG34a2haQ7BdguLjU8AgvzNZgcY%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2012%20Jun%202025%2002%3A01%3A53%20GMT&X-Amz-Signature=0174dd8bb10a01b7509282b3d1b38d5ba97e7e44472073c9b2912c74352957ed"
$ws.Range("L3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1749578078283/BTAP7BioCB-camera-video-segment-1749578080783.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOSQKCUFZ%2F20250611%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250611T180153Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEAIaCXVzLXdlc3QtMiJIMEYCIQDhk5Onjj%2Fjx%2BlzJ8o205%2Bv1tsdlBBdthpD%2F4QpzF6bYwIhALNEtuuSmacwYjQvcc0fMlCNBKXYo4J2gNKCXX%2B%2BFrywKuYDCNv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgzUz0HnBnoKLwdoez8qugO3MhgIwL3LclI71xVisID0mBDNs8067PMeYBZe%2Fko42MpEUpV0qkdrO1tB%2B44kCAwuoywLIdW7BVZn1dRlGZoEle6YSkylm5ods18ww5tyy1Opr0i9RNZo%2BWs433CKDcDpNBVuGQ%2FGLm9iJJSabX8z6Ncao2cxw5RY7TBfDD5rs5JyPkyZLqDrqyL2wESG2b%2BMSQDTjjuv5%2B9dZKJlmlS7IbzOKAV7i3L4Nq3oRxfEhqClmurbIC0HYYb%2FIpqOVoeR7snITMw3Fs2FKMJuZogmbEwxShPXmf%2FOT7pD62oVaAXztypnUSALnv9ijYdXa9khUkanDDiiHoqZh7EdHrh9XwYQwi6wHj03C%2FRsQCWs2e3IyKkc0KwPpBnW6qRa3%2FSZ5Wr5Vj2WUKkSD9apEfa%2FxAwj5aTpn3h%2FdkLWTdI5s1ajKq3FXIcogbcrtGHxSnPXo9JKDh53iscfm7ltH%2Bf5E24SkIcBuxL0oDtgUmZPzU8%2FqpGPF3TOi0Vh05rzhZ6hyzR62J7%2Bhnd4PzqLIMjNlJqhjlD6lGlgW4SRix1egzOpUWWPZ6g05OY7J8TcYqnvh%2BzwGPNzIA1BMLeAp8IGOqQBbYFp2vYp3UIZhng%2B7s9p23ghBnFOaP8tnV%2FZYfvpObScfaqtZixTQfb4CohqRC5nEZNktkMpWTNjEiBcvAyB6bYR8b4fiLQGFpaY0xMTYZZhp3%2FrIEyLSyD8Tl0XWURGYalzMic%2Bf4IYUvGCh5lTYLJ3%2FX8ltWZzshJtOT0LY80f0Q6mYJoFgv%2F7677VGhbhlG34a2haQ7BdguLjU8AgvzNZgcY%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2012%20Jun%202025%2002%3A01%3A53%20GMT&X-Amz-Signature=d70da66856d151ff5c14f5f700cae9b39cdd0f09d33e27635e1eaadea3eb778e"
